# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" right after "总计" (and before
#    "2022-Q2"), populated with the Q3 fund-holding detail row.
# 2. Insert a new summary row for "2022-Q3" at the top of the data table on
#    the "总计" sheet, pushing the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# Writing a numeric-looking string via .Value lets Excel auto-convert it to
# a real number. Several columns in these sheets store numeric-looking
# figures as plain text, so force text storage: flip the cell to a text
# number-format before assigning, then paste-special just the *formats*
# from a pristine, never-touched cell back on top so no stray per-cell
# number-format sticks around afterwards.
function Set-TextValue($cell, $text) {
    $sheet = $cell.Worksheet
    $blank = $sheet.Cells.Item($sheet.Rows.Count, $sheet.Columns.Count)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $blank.Copy()
    $cell.PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------------
# Step 1: create the "2022-Q3" sheet right after "总计".
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$q3Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q3Sheet.Name = "2022-Q3"

# Worksheets.Add() shifts the position of every sheet that came after the
# insertion point, so re-resolve "2022-Q2" *after* the insert rather than
# reuse a handle obtained beforehand (stale handles read back blank values).
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# Copy the header row + first data row formatting from the "2022-Q2" sheet
# so the new sheet matches the look of its siblings, then overwrite with
# the Q3 figures. Column A is left out of the header-row copy: it is blank
# in row 1 and copying it would materialise a spurious empty A1 cell.
$q2Sheet.Range("B1:H1").Copy($q3Sheet.Range("B1:H1"))
$q2Sheet.Range("A2:H2").Copy($q3Sheet.Range("A2:H2"))

$q3Sheet.Cells.Item(2, 1).Value = 0
Set-TextValue $q3Sheet.Cells.Item(2, 2) "400032"
$q3Sheet.Cells.Item(2, 3).Value = "东方主题精选混合"
Set-TextValue $q3Sheet.Cells.Item(2, 4) "14.06"
Set-TextValue $q3Sheet.Cells.Item(2, 5) "87.68"
Set-TextValue $q3Sheet.Cells.Item(2, 6) "4.06"
Set-TextValue $q3Sheet.Cells.Item(2, 7) "0.5708"
$q3Sheet.Cells.Item(2, 8).Value = 4

# ---------------------------------------------------------------------------
# Step 2: insert the "2022-Q3" row into the "总计" summary table.
# ---------------------------------------------------------------------------

# Snapshot the existing data rows (2..4) before they get shifted down.
$existing = @()
for ($r = 2; $r -le 4; $r++) {
    $existing += , @(
        $totalSheet.Cells.Item($r, 2).Value2,
        $totalSheet.Cells.Item($r, 3).Value2,
        $totalSheet.Cells.Item($r, 4).Value2
    )
}

# Make room for the new row 5 (copy formatting of row 4 down one row so the
# index-column style is preserved for the newly created row).
$totalSheet.Range("A4").Copy($totalSheet.Range("A5"))

# Re-write rows 3..5 with the snapshotted values (shifted down by one row);
# column A is a 0-based running index, as in the original sheet.
for ($i = 0; $i -lt $existing.Length; $i++) {
    $r = 3 + $i
    $row = $existing[$i]
    $totalSheet.Cells.Item($r, 1).Value = $i + 1
    $totalSheet.Cells.Item($r, 2).Value = $row[0]
    $totalSheet.Cells.Item($r, 3).Value = $row[1]
    $totalSheet.Cells.Item($r, 4).Value = $row[2]
}

# Write the new "2022-Q3" row at row 2.
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 1
$totalSheet.Cells.Item(2, 4).Value = 0.57
